$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header shared strings: volume number and report week dates ---
$ws.Range("A8").Value = "Volume 31   Number  16"
$ws.Range("C9").Value = "Report Covering the Week  4/15/2024  Through  4/21/2024"

# --- Cells that change between numeric and placeholder-text representation ---
# Donor cells (stable / unmodified by this edit) used purely as a formatting
# template so the pasted cell keeps the exact existing style used elsewhere
# in the sheet for that representation (text "0", text "***.*", plain number,
# or percent-style number).
$ws.Range("C22").Copy($ws.Range("C15"))  # -> text "0"
$ws.Range("F22").Copy($ws.Range("D15"))
$ws.Range("D15").Value = 1
$ws.Range("L22").Copy($ws.Range("E15"))
$ws.Range("E15").Value = -100
$ws.Range("C22").Copy($ws.Range("D23"))  # -> text "0"
$ws.Range("N22").Copy($ws.Range("E23"))  # -> text "***.*"
$ws.Range("C22").Copy($ws.Range("C27"))  # -> text "0"
$ws.Range("F22").Copy($ws.Range("D27"))
$ws.Range("D27").Value = 1
$ws.Range("L22").Copy($ws.Range("E27"))
$ws.Range("E27").Value = -100

# --- Plain numeric value updates (style/type unchanged) ---
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 4
$ws.Range("K15").Value = 25
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 150
$ws.Range("F16").Value = 12
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = 9.090909090909
$ws.Range("I16").Value = 56
$ws.Range("J16").Value = 55
$ws.Range("K16").Value = 1.818181818181
$ws.Range("L16").Value = -16.417910447761
$ws.Range("M16").Value = 33.333333333333
$ws.Range("N16").Value = -86.341463414634
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = -33.333333333333
$ws.Range("F17").Value = 11
$ws.Range("G17").Value = 20
$ws.Range("H17").Value = -45
$ws.Range("I17").Value = 40
$ws.Range("J17").Value = 55
$ws.Range("K17").Value = -27.272727272727
$ws.Range("L17").Value = -18.367346938775
$ws.Range("M17").Value = 42.857142857142
$ws.Range("N17").Value = -55.056179775280
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 8
$ws.Range("E18").Value = -62.5
$ws.Range("G18").Value = 24
$ws.Range("H18").Value = -41.666666666666
$ws.Range("I18").Value = 71
$ws.Range("J18").Value = 93
$ws.Range("K18").Value = -23.655913978494
$ws.Range("L18").Value = -6.578947368421
$ws.Range("M18").Value = -21.978021978022
$ws.Range("N18").Value = -91.486810551558
$ws.Range("C19").Value = 25
$ws.Range("D19").Value = 35
$ws.Range("E19").Value = -28.571428571428
$ws.Range("F19").Value = 99
$ws.Range("G19").Value = 127
$ws.Range("H19").Value = -22.047244094488
$ws.Range("I19").Value = 462
$ws.Range("J19").Value = 466
$ws.Range("K19").Value = -0.858369098712
$ws.Range("L19").Value = -0.645161290322
$ws.Range("M19").Value = 25.543478260869
$ws.Range("N19").Value = -59.791122715404
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 8
$ws.Range("E20").Value = -75
$ws.Range("F20").Value = 5
$ws.Range("G20").Value = 14
$ws.Range("H20").Value = -64.285714285714
$ws.Range("I20").Value = 19
$ws.Range("J20").Value = 38
$ws.Range("K20").Value = -50
$ws.Range("L20").Value = -47.222222222222
$ws.Range("M20").Value = 11.764705882352
$ws.Range("N20").Value = -98.224299065420
$ws.Range("C21").Value = 37
$ws.Range("D21").Value = 57
$ws.Range("E21").Value = -35.087719298245
$ws.Range("F21").Value = 143
$ws.Range("G21").Value = 198
$ws.Range("H21").Value = -27.777777777777
$ws.Range("I21").Value = 653
$ws.Range("J21").Value = 713
$ws.Range("K21").Value = -8.415147265077
$ws.Range("L21").Value = -6.312769010043
$ws.Range("M21").Value = 18.511796733212
$ws.Range("N21").Value = -81.657303370786
$ws.Range("D22").Value = 2
$ws.Range("G22").Value = 5
$ws.Range("H22").Value = -80
$ws.Range("J22").Value = 13
$ws.Range("K22").Value = -15.384615384615
$ws.Range("M22").Value = 83.333333333333
$ws.Range("M23").Value = -20
$ws.Range("C24").Value = 60
$ws.Range("D24").Value = 51
$ws.Range("E24").Value = 17.647058823529
$ws.Range("F24").Value = 212
$ws.Range("G24").Value = 218
$ws.Range("H24").Value = -2.752293577981
$ws.Range("I24").Value = 878
$ws.Range("J24").Value = 925
$ws.Range("K24").Value = -5.081081081081
$ws.Range("L24").Value = -20.614828209764
$ws.Range("M24").Value = 92.543859649122
$ws.Range("C25").Value = 49
$ws.Range("D25").Value = 47
$ws.Range("E25").Value = 4.255319148936
$ws.Range("F25").Value = 178
$ws.Range("G25").Value = 177
$ws.Range("H25").Value = 0.564971751412
$ws.Range("I25").Value = 773
$ws.Range("J25").Value = 787
$ws.Range("K25").Value = -1.778907242693
$ws.Range("L25").Value = -24.066797642436
$ws.Range("C26").Value = 8
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = 14.285714285714
$ws.Range("F26").Value = 24
$ws.Range("G26").Value = 28
$ws.Range("H26").Value = -14.285714285714
$ws.Range("I26").Value = 103
$ws.Range("J26").Value = 113
$ws.Range("K26").Value = -8.849557522123
$ws.Range("L26").Value = -8.849557522123
$ws.Range("M26").Value = -0.961538461538
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 7
$ws.Range("K27").Value = -28.571428571428
$ws.Range("C28").Value = 2
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 11
$ws.Range("G28").Value = 9
$ws.Range("H28").Value = 22.222222222222
$ws.Range("I28").Value = 41
$ws.Range("J28").Value = 30
$ws.Range("K28").Value = 36.666666666666
$ws.Range("L28").Value = 57.692307692307
$ws.Range("F31").Value = 3
